# Include first name of each athlete in 1000K race results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = "HOKIA LINTITA"
$ws.Range("B5").Value = "SANDOR BOGI"
$ws.Range("B4").Value = "JAROSLAV PRUCKNER"

$ws.Columns.Item(2).AutoFit() | Out-Null

$ws.Range("B6").Select() | Out-Null
